$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Duplicate A3 ("EV002") into C3, extending the sheet's used range to column C.
# Use .Text (not .Value, which round-trips through a COM Variant wrapper here)
# so the literal string is written rather than a boxed object representation.
$a3val = $ws.Range("A3").Text
$ws.Range("C3").Value = $a3val
